$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650477876173526"
$ws1.Range("B2").Value = "go_stims-1650477876135491.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778761565275.csv"
$ws1.Range("B4").Value = "go_stims-16504778761574938.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778761725276.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778773623347"
$ws2.Range("B2").Value = "ZB-match_3-16504778761954908.csv"
$ws2.Range("B3").Value = "ZB-match_9-16504778764075222.csv"
$ws2.Range("B4").Value = "TB-16504778769953644.csv"
$ws2.Range("B5").Value = "OB-16504778768455245.csv"
$ws2.Range("B6").Value = "ZB-match_0-16504778765004919.csv"
$ws2.Range("B7").Value = "TB-1650477877011329.csv"
$ws2.Range("B8").Value = "OB-16504778768885252.csv"
$ws2.Range("B9").Value = "TB-16504778773443651.csv"
$ws2.Range("B10").Value = "OB-16504778766674917.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650477877363331"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778774103315"
$ws4.Range("B2").Value = "MM_stims-16504778773773649.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778773653316.csv"
$ws4.Range("B4").Value = "MM_stims-16504778773933637.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778773773649.csv"
$ws4.Range("B6").Value = "MM_stims-16504778774093642.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778773943307.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778774733908"
$ws5.Range("B2").Value = "SAT_stims-1650477877426369.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778774574287.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778774143324.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477877441391.csv"
